$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 45 new rows (102-146) following the existing repeating pattern:
#   A cycles 10002..10010 (9 values), B cycles 10021..10029 (9 values),
#   C increments by 1 each row starting at 3000121, and D-H stay constant.
$startRow = 102
$startC = 3000121

for ($i = 0; $i -lt 45; $i++) {
    $row = $startRow + $i
    $a = 10002 + ($i % 9)
    $b = 10021 + ($i % 9)
    $c = $startC + $i

    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin()"
    $ws.Cells.Item($row, 7).Value = "now()"
    $ws.Cells.Item($row, 8).Value = "now()"
}

# Mirror the selection left behind by the author (cursor parked on the
# row right after the new data, selecting down to the end of the sheet).
$ws.Range("A147:XFD1048576").Select()

# Set up page printing as in the committed version.
$ws.PageSetup.Orientation = 1
